# Add the next day's "Daily User Impact Status" row (row 23) below the
# existing data, matching the same layout/formatting as prior rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date column: store as a date serial and apply the "d-mmm-yy" number
# format (numFmtId 15), matching the new style added for this row.
$ws.Range("A23").Value = 45967
$ws.Range("A23").NumberFormat = "d-mmm-yy"

# Remaining metric columns use the sheet's normal (unformatted) style.
$ws.Range("B23").Value = 5597
$ws.Range("C23").Value = 4327
$ws.Range("D23").Value = 3993
$ws.Range("E23").Value = 258
$ws.Range("F23").Value = 46
$ws.Range("G23").Value = 26
$ws.Range("H23").Value = 4
$ws.Range("I23").Value = 0

# Move the active selection onto the newly added row, as happens after
# entering data in the next blank row.
[void]$ws.Range("A23:I23").Select()
